# Auto-generated COM-interop script applying the scraped-sheet value updates.
# For each touched cell: set new numeric value, or clear the cell entirely
# when the diff shows the <c> element being removed outright.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 457.7143
$ws.Range("I33").Value = 292.41666
$ws.Range("K33").Value = 292.41666
$ws.Range("M33").Value = -63.41665999999998
$ws.Range("H54").Value = 40838.8
$ws.Range("H80").Value = 6761.476
$ws.Range("I80").Value = 1199.4166
$ws.Range("K80").Value = 3598.2498
$ws.Range("M80").Value = -2600.2498
$ws.Range("H83").Value = 6761.476
$ws.Range("I83").Value = 1199.4166
$ws.Range("K83").Value = 10794.7494
$ws.Range("M83").Value = -5802.749400000001
$ws.Range("H88").Value = 1503267.2
$ws.Range("J88").Value = 3534.6667
$ws.Range("L88").Value = 3534.6667
$ws.Range("N88").Value = -4346.6667
$ws.Range("H91").Value = 1503267.2
$ws.Range("J91").Value = 3534.6667
$ws.Range("L91").Value = 3534.6667
$ws.Range("N91").Value = -6342.6667
$ws.Range("H100").Value = 2658.5334
$ws.Range("I100").Value = 2374.5454
$ws.Range("K100").Value = 2374.5454
$ws.Range("M100").Value = -1833.5454
$ws.Range("H116").Value = 8136.273
$ws.Range("I116").Value = 7032
$ws.Range("J116").Value = 10068.75
$ws.Range("K116").Value = 7032
$ws.Range("L116").Value = 10068.75
$ws.Range("M116").Value = -3590
$ws.Range("N116").Value = -16952.75
$ws.Range("H137").Value = 280694.97
$ws.Range("I137").Value = 1806.4348
$ws.Range("J137").Value = 501882.4
$ws.Range("K137").Value = 5419.3044
$ws.Range("L137").Value = 1505647.2
$ws.Range("M137").Value = -2869.3044
$ws.Range("N137").Value = -1510747.2
$ws.Range("H138").Value = 1957.4791
$ws.Range("J138").Value = 2449.48
$ws.Range("L138").Value = 7348.440000000001
$ws.Range("N138").Value = -17628.44

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1196.6666
$ws.Range("J2").Value = 1422.5
$ws.Range("L2").Value = 1422.5
$ws.Range("N2").Value = -1648.5
$ws.Range("H45").Value = 11369429
$ws.Range("I45").Value = 2670.5
$ws.Range("J45").Value = 25009540
$ws.Range("K45").Value = 2670.5
$ws.Range("L45").Value = 25009540
$ws.Range("M45").Value = -2293.5
$ws.Range("N45").Value = -25010294
$ws.Range("H60").Value = 5799
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 5799
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 5799
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -7265
$ws.Range("H74").Value = 29657.055
$ws.Range("I74").Value = 42851.25
$ws.Range("J74").Value = 3268.6667
$ws.Range("K74").Value = 42851.25
$ws.Range("L74").Value = 3268.6667
$ws.Range("M74").Value = -41977.25
$ws.Range("N74").Value = -5016.6667
$ws.Range("H77").Value = 29657.055
$ws.Range("I77").Value = 42851.25
$ws.Range("J77").Value = 3268.6667
$ws.Range("K77").Value = 214256.25
$ws.Range("L77").Value = 16343.3335
$ws.Range("M77").Value = -209888.25
$ws.Range("N77").Value = -25079.3335
$ws.Range("H110").Value = 397.76
$ws.Range("I110").Value = 437.2381
$ws.Range("K110").Value = 437.2381
$ws.Range("M110").Value = 1607.7619
$ws.Range("H116").Value = 1196.6666
$ws.Range("J116").Value = 1422.5
$ws.Range("L116").Value = 1422.5
$ws.Range("N116").Value = -6010.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1196.6666
$ws.Range("J3").Value = 1422.5
$ws.Range("L3").Value = 1422.5
$ws.Range("N3").Value = -1650.5
$ws.Range("H22").Value = 25974690
$ws.Range("I22").Value = 25974690
$ws.Range("K22").Value = 25974690
$ws.Range("M22").Value = -25974517
$ws.Range("H54").Value = 2950.4443
$ws.Range("I54").Value = 1444.25
$ws.Range("K54").Value = 1444.25
$ws.Range("M54").Value = -960.25
$ws.Range("H99").Value = 441773.38
$ws.Range("I99").Value = 72939.5
$ws.Range("J99").Value = 911198.25
$ws.Range("K99").Value = 72939.5
$ws.Range("L99").Value = 911198.25
$ws.Range("M99").Value = -71441.5
$ws.Range("N99").Value = -914194.25
$ws.Range("H134").Value = 5600
$ws.Range("I134").Value = 1360
$ws.Range("K134").Value = 4080
$ws.Range("M134").Value = -1545

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1831.2941
$ws.Range("I16").Value = 1653.6
$ws.Range("J16").Value = 2085.1428
$ws.Range("K16").Value = 1653.6
$ws.Range("L16").Value = 2085.1428
$ws.Range("M16").Value = -1366.6
$ws.Range("N16").Value = -2659.1428
$ws.Range("H55").Value = 15024.333
$ws.Range("I55").Value = 5073
$ws.Range("J55").Value = 20000
$ws.Range("K55").Value = 5073
$ws.Range("L55").Value = 20000
$ws.Range("M55").Value = -4758
$ws.Range("N55").Value = -20630
$ws.Range("H97").Value = 45000
$ws.Range("J97").Value = 45000
$ws.Range("L97").Value = 45000
$ws.Range("N97").Value = -46982
$ws.Range("H105").Value = 2430.3076
$ws.Range("I105").Value = 849.5
$ws.Range("K105").Value = 849.5
$ws.Range("M105").Value = 897.5
$ws.Range("H107").Value = 1173.5714
$ws.Range("I107").Value = 1237.8
$ws.Range("J107").Value = 1013
$ws.Range("K107").Value = 1237.8
$ws.Range("L107").Value = 1013
$ws.Range("M107").Value = 682.2
$ws.Range("N107").Value = -4853
$ws.Range("H113").Value = 1831.2941
$ws.Range("I113").Value = 1653.6
$ws.Range("J113").Value = 2085.1428
$ws.Range("K113").Value = 1653.6
$ws.Range("L113").Value = 2085.1428
$ws.Range("M113").Value = 516.4000000000001
$ws.Range("N113").Value = -6425.1428
$ws.Range("H122").Value = 2402.8635
$ws.Range("J122").Value = 2699
$ws.Range("L122").Value = 8097
$ws.Range("N122").Value = -12997
$ws.Range("H134").Value = 43280.168
$ws.Range("I134").Value = 1274.5834
$ws.Range("K134").Value = 3823.7502
$ws.Range("M134").Value = -1288.7502

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 1804.5151
$ws.Range("I140").Value = 961.3182
$ws.Range("J140").Value = 3490.9092
$ws.Range("K140").Value = 2883.9546
$ws.Range("L140").Value = 10472.7276
$ws.Range("M140").Value = 2296.0454
$ws.Range("N140").Value = -20832.7276

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 1918
$ws.Range("I55").Value = 1433.3334
$ws.Range("J55").Value = 2645
$ws.Range("K55").Value = 1433.3334
$ws.Range("L55").Value = 2645
$ws.Range("M55").Value = -1106.3334
$ws.Range("N55").Value = -3299
$ws.Range("H70").Value = 118411.375
$ws.Range("I70").Value = 7055.4
$ws.Range("J70").Value = 304004.66
$ws.Range("K70").Value = 7055.4
$ws.Range("L70").Value = 304004.66
$ws.Range("M70").Value = -6785.4
$ws.Range("N70").Value = -304544.66
$ws.Range("H73").Value = 118411.375
$ws.Range("I73").Value = 7055.4
$ws.Range("J73").Value = 304004.66
$ws.Range("K73").Value = 7055.4
$ws.Range("L73").Value = 304004.66
$ws.Range("M73").Value = -6119.4
$ws.Range("N73").Value = -305876.66
$ws.Range("H107").Value = 775.4643
$ws.Range("I107").Value = 578.4167
$ws.Range("K107").Value = 578.4167
$ws.Range("M107").Value = 1341.5833
$ws.Range("H123").Value = 25696.715
$ws.Range("J123").Value = 25696.715
$ws.Range("L123").Value = 25696.715
$ws.Range("N123").Value = -30596.715
$ws.Range("H126").Value = 6709.8
$ws.Range("I126").Value = 2460
$ws.Range("K126").Value = 7380
$ws.Range("M126").Value = -4910
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("M139").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5769.2
$ws.Range("I22").Value = 987.5
$ws.Range("J22").Value = 8957
$ws.Range("K22").Value = 987.5
$ws.Range("L22").Value = 8957
$ws.Range("M22").Value = -692.5
$ws.Range("N22").Value = -9547
$ws.Range("H27").Value = 5769.2
$ws.Range("I27").Value = 987.5
$ws.Range("J27").Value = 8957
$ws.Range("K27").Value = 987.5
$ws.Range("L27").Value = 8957
$ws.Range("M27").Value = -880.5
$ws.Range("N27").Value = -9171
$ws.Range("H46").Value = 10595.846
$ws.Range("I46").Value = 16363.857
$ws.Range("J46").Value = 3866.5
$ws.Range("K46").Value = 16363.857
$ws.Range("L46").Value = 3866.5
$ws.Range("M46").Value = -16175.857
$ws.Range("N46").Value = -4242.5
$ws.Range("H61").Value = 2098.85
$ws.Range("I61").Value = 1857.9286
$ws.Range("K61").Value = 1857.9286
$ws.Range("M61").Value = -1655.9286
$ws.Range("H68").Value = 4428.857
$ws.Range("I68").Value = 4720.4
$ws.Range("J68").Value = 3700
$ws.Range("K68").Value = 4720.4
$ws.Range("L68").Value = 3700
$ws.Range("M68").Value = -3971.4
$ws.Range("N68").Value = -5198
$ws.Range("H71").Value = 4428.857
$ws.Range("I71").Value = 4720.4
$ws.Range("J71").Value = 3700
$ws.Range("K71").Value = 23602
$ws.Range("L71").Value = 18500
$ws.Range("M71").Value = -19858
$ws.Range("N71").Value = -25988
$ws.Range("H80").Value = 29998
$ws.Range("J80").Value = 29998
$ws.Range("L80").Value = 29998
$ws.Range("N80").Value = -32244
$ws.Range("H83").Value = 29998
$ws.Range("J83").Value = 29998
$ws.Range("L83").Value = 89994
$ws.Range("N83").Value = -101226
$ws.Range("H113").Value = 2098.85
$ws.Range("I113").Value = 1857.9286
$ws.Range("K113").Value = 1857.9286
$ws.Range("M113").Value = 312.0714

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 5862674.5
$ws.Range("I96").Value = 21100
$ws.Range("J96").Value = 13164642
$ws.Range("K96").Value = 21100
$ws.Range("L96").Value = 13164642
$ws.Range("M96").Value = -19727
$ws.Range("N96").Value = -13167388
$ws.Range("H122").Value = 2219.5557
$ws.Range("I122").Value = 1404
$ws.Range("K122").Value = 4212
$ws.Range("M122").Value = -1762
$ws.Range("H124").Value = 84999.164
$ws.Range("J124").Value = 84999.164
$ws.Range("L124").Value = 84999.164
$ws.Range("N124").Value = -94819.164
